# min-eV PCBv1.2.2: correct PCBway BOM (OD90 / OD135 resistors were swapped)
# and move the saved cursor/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The BOM rows for the two resistors (row 10 = OD90, row 11 = OD135) had their
# Manufacturer Part Number (col E) and Description/Value (col F) swapped by
# mistake. Swap the two rows' values back so each designator group gets the
# correct part.
$e10 = $ws.Range("E10").Value()
$f10 = $ws.Range("F10").Value()
$e11 = $ws.Range("E11").Value()
$f11 = $ws.Range("F11").Value()

$ws.Range("E10").Value = $e11
$ws.Range("F10").Value = $f11
$ws.Range("E11").Value = $e10
$ws.Range("F11").Value = $f10

# Move the active selection to where the author left off editing.
$ws.Range("F13").Select()
